# Weekly update: insert the newest week's price row at the top of the
# data block (row 7), pushing all the existing data rows down by one.
# This mirrors the source system behaviour where the freshest
# observation is prepended above the previously-newest rows while the
# very first few summary rows (2-6) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at row 7; everything from the old row 7 onward
# (through the old row 117) shifts down to rows 8-118, growing the
# used range from A1:R117 to A1:R118.
$ws.Rows("7:7").Insert()

# Populate the new row 7 with this week's observation.
$ws.Cells.Item(7, 1).Value  = 1
$ws.Cells.Item(7, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value  = 45083
$ws.Cells.Item(7, 5).Value  = 15
$ws.Cells.Item(7, 6).Value  = 100112040
$ws.Cells.Item(7, 7).Value  = "Cilantro"
$ws.Cells.Item(7, 8).Value  = "Sin especificar"
$ws.Cells.Item(7, 9).Value  = "Primera"
$ws.Cells.Item(7, 10).Value = 350
$ws.Cells.Item(7, 11).Value = 2400
$ws.Cells.Item(7, 12).Value = 2500
$ws.Cells.Item(7, 13).Value = 2443
$ws.Cells.Item(7, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 1222
$ws.Cells.Item(7, 17).Value = 2
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# Give the new date cell the same date number-format as the rest of
# column D (style carries over automatically from Insert, but set it
# explicitly too so the value renders as a date rather than a serial).
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
